$d = $word.ActiveDocument

# The document has one section; each section carries a "default" header/
# footer (index 1) and a "first page" header/footer (index 2), since
# Different First Page is turned on (w:titlePg).
#
#   Headers.Item(1) -> default header  (BTEC logo, currently "image1.jpg")
#   Headers.Item(2) -> first-page header (BTEC logo, currently "image1.jpg")
#   Footers.Item(1) -> default footer  (Pearson logo, currently "image2.png")
#   Footers.Item(2) -> first-page footer (Pearson logo, currently "image2.png")
#
# Rename each embedded picture: the two BTEC logos from image1.jpg to
# image2.jpg, and the two Pearson logos from image2.png to image1.png.

$sec = $d.Sections.Item(1)

$headerDefault = $sec.Headers.Item(1)
if ($headerDefault.Exists -and $headerDefault.Range.InlineShapes.Count -ge 1) {
    $headerDefault.Range.InlineShapes.Item(1).Name = "image2.jpg"
}

$headerFirst = $sec.Headers.Item(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    $headerFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"
}

$footerDefault = $sec.Footers.Item(1)
if ($footerDefault.Exists -and $footerDefault.Range.InlineShapes.Count -ge 1) {
    $footerDefault.Range.InlineShapes.Item(1).Name = "image1.png"
}

$footerFirst = $sec.Footers.Item(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    $footerFirst.Range.InlineShapes.Item(1).Name = "image1.png"
}

Write-Output "Renamed header/footer logo images."
